# The deck's theme colour palette changes from the "Integral" scheme to
# the "Office Theme" scheme (the slide master's theme and the notes
# master's theme effectively swap colour palettes).
#
# PowerPoint's object model exposes the slide master's 12 theme colour
# slots via SlideMaster.Theme.ThemeColorScheme.Colors(i).RGB, mapping
# 1:1 onto <a:clrScheme>'s children in document order: dk1, lt1, dk2,
# lt2, accent1-6, hlink, folHlink. COM's RGB long packs bytes as
# 0xBBGGRR, so each target hex colour RRGGBB is written as 0xBBGGRR.

$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$tcs = $m.Theme.ThemeColorScheme

$tcs.Colors(1).RGB  = 0x000000   # dk1       -> 000000
$tcs.Colors(2).RGB  = 0xFFFFFF   # lt1       -> FFFFFF
$tcs.Colors(3).RGB  = 0x6A5444   # dk2       -> 44546A
$tcs.Colors(4).RGB  = 0xE6E6E7   # lt2       -> E7E6E6
$tcs.Colors(5).RGB  = 0xD59B5B   # accent1   -> 5B9BD5
$tcs.Colors(6).RGB  = 0x317DED   # accent2   -> ED7D31
$tcs.Colors(7).RGB  = 0xA5A5A5   # accent3   -> A5A5A5
$tcs.Colors(8).RGB  = 0x00C0FF   # accent4   -> FFC000
$tcs.Colors(9).RGB  = 0xC47244   # accent5   -> 4472C4
$tcs.Colors(10).RGB = 0x47AD70   # accent6   -> 70AD47
$tcs.Colors(11).RGB = 0xC16305   # hlink     -> 0563C1
$tcs.Colors(12).RGB = 0x724F95   # folHlink  -> 954F72
